$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dStyle = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "19.962.34"
$ws.Range("D2").Style = $dStyle
$eStyle = $ws.Range("E2").Style
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -3.18%  "
$ws.Range("E2").Style = $eStyle

$dStyle = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.416.78"
$ws.Range("D3").Style = $dStyle
$eStyle = $ws.Range("E3").Style
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.33%  "
$ws.Range("E3").Style = $eStyle

$dStyle = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = $dStyle
$eStyle = $ws.Range("E4").Style
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.70%  "
$ws.Range("E4").Style = $eStyle

$dStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.001"
$ws.Range("D5").Style = $dStyle
$eStyle = $ws.Range("E5").Style
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.24%  "
$ws.Range("E5").Style = $eStyle

$dStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "276.43"
$ws.Range("D6").Style = $dStyle
$eStyle = $ws.Range("E6").Style
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.73%  "
$ws.Range("E6").Style = $eStyle

$dStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3702"
$ws.Range("D7").Style = $dStyle
$eStyle = $ws.Range("E7").Style
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.72%  "
$ws.Range("E7").Style = $eStyle

$dStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3110"
$ws.Range("D8").Style = $dStyle
$eStyle = $ws.Range("E8").Style
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.80%  "
$ws.Range("E8").Style = $eStyle

$dStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.90"
$ws.Range("D9").Style = $dStyle
$eStyle = $ws.Range("E9").Style
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.96%  "
$ws.Range("E9").Style = $eStyle

$dStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.038"
$ws.Range("D10").Style = $dStyle
$eStyle = $ws.Range("E10").Style
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.87%  "
$ws.Range("E10").Style = $eStyle

$dStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06539"
$ws.Range("D11").Style = $dStyle
$eStyle = $ws.Range("E11").Style
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.40%  "
$ws.Range("E11").Style = $eStyle

$dStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("D12").Style = $dStyle
$eStyle = $ws.Range("E12").Style
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.73%  "
$ws.Range("E12").Style = $eStyle

$dStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.479"
$ws.Range("D13").Style = $dStyle
$eStyle = $ws.Range("E13").Style
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.97%  "
$ws.Range("E13").Style = $eStyle

$dStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.68"
$ws.Range("D14").Style = $dStyle
$eStyle = $ws.Range("E14").Style
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.78%  "
$ws.Range("E14").Style = $eStyle

$dStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.211"
$ws.Range("D15").Style = $dStyle
$eStyle = $ws.Range("E15").Style
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.44%  "
$ws.Range("E15").Style = $eStyle

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$dStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.420.55"
$ws.Range("D16").Style = $dStyle
$eStyle = $ws.Range("E16").Style
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.39%  "
$ws.Range("E16").Style = $eStyle

$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$dStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001022"
$ws.Range("D17").Style = $dStyle
$eStyle = $ws.Range("E17").Style
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("E17").Style = $eStyle

$dStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.05704"
$ws.Range("D18").Style = $dStyle
$eStyle = $ws.Range("E18").Style
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -11.59%  "
$ws.Range("E18").Style = $eStyle

$dStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("D19").Style = $dStyle
$eStyle = $ws.Range("E19").Style
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.21%  "
$ws.Range("E19").Style = $eStyle

$dStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.16"
$ws.Range("D20").Style = $dStyle
$eStyle = $ws.Range("E20").Style
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -8.96%  "
$ws.Range("E20").Style = $eStyle

$dStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.613"
$ws.Range("D21").Style = $dStyle
$eStyle = $ws.Range("E21").Style
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.92%  "
$ws.Range("E21").Style = $eStyle

$dStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.79"
$ws.Range("D22").Style = $dStyle
$eStyle = $ws.Range("E22").Style
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.99%  "
$ws.Range("E22").Style = $eStyle

$dStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.99"
$ws.Range("D23").Style = $dStyle
$eStyle = $ws.Range("E23").Style
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.48%  "
$ws.Range("E23").Style = $eStyle

$dStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.239"
$ws.Range("D24").Style = $dStyle
$eStyle = $ws.Range("E24").Style
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -4.24%  "
$ws.Range("E24").Style = $eStyle

$dStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "19.998.54"
$ws.Range("D25").Style = $dStyle
$eStyle = $ws.Range("E25").Style
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.88%  "
$ws.Range("E25").Style = $eStyle

$dStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.274"
$ws.Range("D26").Style = $dStyle
$eStyle = $ws.Range("E26").Style
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.78%  "
$ws.Range("E26").Style = $eStyle

$dStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "133.21"
$ws.Range("D27").Style = $dStyle
$eStyle = $ws.Range("E27").Style
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -6.96%  "
$ws.Range("E27").Style = $eStyle

$dStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.31"
$ws.Range("D28").Style = $dStyle
$eStyle = $ws.Range("E28").Style
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.31%  "
$ws.Range("E28").Style = $eStyle

$dStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.580.35"
$ws.Range("D29").Style = $dStyle
$eStyle = $ws.Range("E29").Style
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.82%  "
$ws.Range("E29").Style = $eStyle

$dStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "110.12"
$ws.Range("D30").Style = $dStyle
$eStyle = $ws.Range("E30").Style
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.37%  "
$ws.Range("E30").Style = $eStyle

$dStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.915"
$ws.Range("D31").Style = $dStyle
$eStyle = $ws.Range("E31").Style
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +13.62%  "
$ws.Range("E31").Style = $eStyle

$dStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.222"
$ws.Range("D32").Style = $dStyle
$eStyle = $ws.Range("E32").Style
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -5.70%  "
$ws.Range("E32").Style = $eStyle

$dStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8144"
$ws.Range("D33").Style = $dStyle
$eStyle = $ws.Range("E33").Style
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -11.60%  "
$ws.Range("E33").Style = $eStyle

$dStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07775"
$ws.Range("D34").Style = $dStyle
$eStyle = $ws.Range("E34").Style
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.26%  "
$ws.Range("E34").Style = $eStyle

$dStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.484"
$ws.Range("D35").Style = $dStyle
$eStyle = $ws.Range("E35").Style
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.49%  "
$ws.Range("E35").Style = $eStyle

$dStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.904"
$ws.Range("D36").Style = $dStyle
$eStyle = $ws.Range("E36").Style
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.11%  "
$ws.Range("E36").Style = $eStyle

$dStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.162"
$ws.Range("D37").Style = $dStyle
$eStyle = $ws.Range("E37").Style
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.17%  "
$ws.Range("E37").Style = $eStyle

$dStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05859"
$ws.Range("D38").Style = $dStyle
$eStyle = $ws.Range("E38").Style
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +4.47%  "
$ws.Range("E38").Style = $eStyle

$dStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.001"
$ws.Range("D39").Style = $dStyle
$eStyle = $ws.Range("E39").Style
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.25%  "
$ws.Range("E39").Style = $eStyle

$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$dStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02052"
$ws.Range("D40").Style = $dStyle
$eStyle = $ws.Range("E40").Style
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("E40").Style = $eStyle

$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$dStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.51"
$ws.Range("D41").Style = $dStyle
$eStyle = $ws.Range("E41").Style
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.91%  "
$ws.Range("E41").Style = $eStyle

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$dStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.109"
$ws.Range("D42").Style = $dStyle
$eStyle = $ws.Range("E42").Style
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.91%  "
$ws.Range("E42").Style = $eStyle

$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$dStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1879"
$ws.Range("D43").Style = $dStyle
$eStyle = $ws.Range("E43").Style
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.79%  "
$ws.Range("E43").Style = $eStyle

$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$dStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5319"
$ws.Range("D44").Style = $dStyle
$eStyle = $ws.Range("E44").Style
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.85%  "
$ws.Range("E44").Style = $eStyle

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$dStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.36"
$ws.Range("D45").Style = $dStyle
$eStyle = $ws.Range("E45").Style
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.17%  "
$ws.Range("E45").Style = $eStyle

$dStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.536"
$ws.Range("D46").Style = $dStyle
$eStyle = $ws.Range("E46").Style
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.15%  "
$ws.Range("E46").Style = $eStyle

$dStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "116.84"
$ws.Range("D47").Style = $dStyle
$eStyle = $ws.Range("E47").Style
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +6.88%  "
$ws.Range("E47").Style = $eStyle

$dStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5205"
$ws.Range("D48").Style = $dStyle
$eStyle = $ws.Range("E48").Style
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.59%  "
$ws.Range("E48").Style = $eStyle

$dStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.771"
$ws.Range("D49").Style = $dStyle
$eStyle = $ws.Range("E49").Style
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.57%  "
$ws.Range("E49").Style = $eStyle

$dStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.038"
$ws.Range("D50").Style = $dStyle
$eStyle = $ws.Range("E50").Style
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.33%  "
$ws.Range("E50").Style = $eStyle

$dStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.001"
$ws.Range("D51").Style = $dStyle
$eStyle = $ws.Range("E51").Style
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.21%  "
$ws.Range("E51").Style = $eStyle
